$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1205.1904
$ws.Range("I6").Value = 224.6
$ws.Range("J6").Value = 1511.625
$ws.Range("K6").Value = 673.8
$ws.Range("L6").Value = 4534.875
$ws.Range("M6").Value = -561.8
$ws.Range("N6").Value = -4758.875
$ws.Range("H9").Value = 150
$ws.Range("I9").Value = 150
$ws.Range("K9").Value = 150
$ws.Range("M9").Value = 19
$ws.Range("H12").Value = 522.3333
$ws.Range("I12").Value = 347
$ws.Range("J12").Value = 873
$ws.Range("K12").Value = 347
$ws.Range("L12").Value = 873
$ws.Range("M12").Value = -177
$ws.Range("N12").Value = -1213
$ws.Range("H17").Value = 1393.6086
$ws.Range("J17").Value = 1383.6571
$ws.Range("L17").Value = 4150.971299999999
$ws.Range("N17").Value = -4486.971299999999
$ws.Range("H21").Value = 47508.5
$ws.Range("I21").Value = 50017
$ws.Range("J21").Value = 45000
$ws.Range("K21").Value = 50017
$ws.Range("L21").Value = 45000
$ws.Range("M21").Value = -49549
$ws.Range("N21").Value = -45936
$ws.Range("H23").Value = 47508.5
$ws.Range("I23").Value = 50017
$ws.Range("J23").Value = 45000
$ws.Range("K23").Value = 50017
$ws.Range("L23").Value = 45000
$ws.Range("M23").Value = -49783
$ws.Range("N23").Value = -45468
$ws.Range("H38").Value = 2028.2
$ws.Range("J38").Value = 10000
$ws.Range("L38").Value = 30000
$ws.Range("N38").Value = -30744
$ws.Range("H39").Value = 480.875
$ws.Range("I39").Value = 224.5
$ws.Range("J39").Value = 1250
$ws.Range("K39").Value = 673.5
$ws.Range("L39").Value = 3750
$ws.Range("M39").Value = -377.5
$ws.Range("N39").Value = -4342
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H58").Value = 4305
$ws.Range("J58").Value = 5166.6665
$ws.Range("L58").Value = 15499.9995
$ws.Range("N58").Value = -15799.9995
$ws.Range("H62").Value = 250007230
$ws.Range("I62").Value = 333339650
$ws.Range("K62").Value = 333339650
$ws.Range("M62").Value = -333339026
$ws.Range("H65").Value = 250007230
$ws.Range("I65").Value = 333339650
$ws.Range("K65").Value = 1666698250
$ws.Range("M65").Value = -1666695130
$ws.Range("H88").Value = 1136.909
$ws.Range("I88").Value = 1040.6666
$ws.Range("J88").Value = 1252.4
$ws.Range("K88").Value = 1040.6666
$ws.Range("L88").Value = 1252.4
$ws.Range("M88").Value = -634.6666
$ws.Range("N88").Value = -2064.4
$ws.Range("H91").Value = 1136.909
$ws.Range("I91").Value = 1040.6666
$ws.Range("J91").Value = 1252.4
$ws.Range("K91").Value = 1040.6666
$ws.Range("L91").Value = 1252.4
$ws.Range("M91").Value = 363.3334
$ws.Range("N91").Value = -4060.4
$ws.Range("H98").Value = 42804.285
$ws.Range("I98").Value = 56016.5
$ws.Range("K98").Value = 56016.5
$ws.Range("M98").Value = -54518.5
$ws.Range("H106").Value = 10302.529
$ws.Range("I106").Value = 11012.333
$ws.Range("K106").Value = 11012.333
$ws.Range("M106").Value = -10381.333
$ws.Range("H111").Value = 2836.2307
$ws.Range("J111").Value = 2908.2
$ws.Range("L111").Value = 8724.599999999999
$ws.Range("N111").Value = -14858.6
$ws.Range("H116").Value = 5696097
$ws.Range("I116").Value = 6265056.5
$ws.Range("K116").Value = 6265056.5
$ws.Range("M116").Value = -6261614.5
$ws.Range("H122").Value = 42804.285
$ws.Range("I122").Value = 56016.5
$ws.Range("K122").Value = 168049.5
$ws.Range("M122").Value = -165599.5
$ws.Range("H125").Value = 5123.2
$ws.Range("J125").Value = 2227.3333
$ws.Range("L125").Value = 20045.9997
$ws.Range("N125").Value = -24965.9997
$ws.Range("H132").Value = 1669141.1
$ws.Range("I132").Value = 2377.8867
$ws.Range("K132").Value = 7133.6601
$ws.Range("M132").Value = -4603.6601
$ws.Range("H137").Value = 718699.6
$ws.Range("J137").Value = 4999.6665
$ws.Range("L137").Value = 14998.9995
$ws.Range("N137").Value = -20098.9995
$ws.Range("H138").Value = 372235.12
$ws.Range("I138").Value = 1885233.1
$ws.Range("J138").Value = 7028.6895
$ws.Range("K138").Value = 5655699.300000001
$ws.Range("L138").Value = 21086.0685
$ws.Range("M138").Value = -5650559.300000001
$ws.Range("N138").Value = -31366.0685
$ws.Range("H141").Value = 8372.467000000001
$ws.Range("I141").Value = 8613.429
$ws.Range("J141").Value = 4999
$ws.Range("K141").Value = 25840.287
$ws.Range("L141").Value = 14997
$ws.Range("M141").Value = -20660.287
$ws.Range("N141").Value = -25357

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18324.361
$ws.Range("I32").Value = 18324.361
$ws.Range("K32").Value = 18324.361
$ws.Range("M32").Value = -18037.361
$ws.Range("H43").Value = 21750.4
$ws.Range("J43").Value = 21750.4
$ws.Range("L43").Value = 21750.4
$ws.Range("N43").Value = -22376.4
$ws.Range("H45").Value = 92471.44
$ws.Range("I45").Value = 122609.234
$ws.Range("K45").Value = 122609.234
$ws.Range("M45").Value = -122232.234
$ws.Range("H74").Value = 3869.5483
$ws.Range("I74").Value = 15573.444
$ws.Range("J74").Value = 1882.0944
$ws.Range("K74").Value = 15573.444
$ws.Range("L74").Value = 1882.0944
$ws.Range("M74").Value = -14699.444
$ws.Range("N74").Value = -3630.0944
$ws.Range("H77").Value = 3869.5483
$ws.Range("I77").Value = 15573.444
$ws.Range("J77").Value = 1882.0944
$ws.Range("K77").Value = 77867.22
$ws.Range("L77").Value = 9410.472
$ws.Range("M77").Value = -73499.22
$ws.Range("N77").Value = -18146.472
$ws.Range("H110").Value = 2052.4443
$ws.Range("J110").Value = 2283.1428
$ws.Range("L110").Value = 2283.1428
$ws.Range("N110").Value = -6373.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 908.7273
$ws.Range("I7").Value = 571.2857
$ws.Range("J7").Value = 1499.25
$ws.Range("K7").Value = 571.2857
$ws.Range("L7").Value = 1499.25
$ws.Range("M7").Value = -458.2857
$ws.Range("N7").Value = -1725.25
$ws.Range("H22").Value = 709.5833
$ws.Range("J22").Value = 1133.3334
$ws.Range("L22").Value = 1133.3334
$ws.Range("N22").Value = -1833.3334
$ws.Range("H94").Value = 1338.3684
$ws.Range("I94").Value = 303
$ws.Range("K94").Value = 303
$ws.Range("M94").Value = 148
$ws.Range("H99").Value = 6948951.5
$ws.Range("I99").Value = 17858856
$ws.Range("K99").Value = 17858856
$ws.Range("M99").Value = -17857358
$ws.Range("H106").Value = 32390.334
$ws.Range("J106").Value = 32390.334
$ws.Range("L106").Value = 32390.334
$ws.Range("N106").Value = -34914.334
$ws.Range("H122").Value = 11729.846
$ws.Range("I122").Value = 13026.272
$ws.Range("J122").Value = 4599.5
$ws.Range("K122").Value = 39078.81600000001
$ws.Range("L122").Value = 13798.5
$ws.Range("M122").Value = -36628.81600000001
$ws.Range("N122").Value = -18698.5
$ws.Range("H126").Value = 6948951.5
$ws.Range("I126").Value = 17858856
$ws.Range("K126").Value = 53576568
$ws.Range("M126").Value = -53574098
$ws.Range("H132").Value = 5290.778
$ws.Range("I132").Value = 5652.1665
$ws.Range("K132").Value = 16956.4995
$ws.Range("M132").Value = -14426.4995
$ws.Range("H134").Value = 2361.95
$ws.Range("I134").Value = 2153.1875
$ws.Range("K134").Value = 6459.5625
$ws.Range("M134").Value = -3924.5625
$ws.Range("H141").Value = 553205.4399999999
$ws.Range("J141").Value = 579864.5600000001
$ws.Range("L141").Value = 579864.5600000001
$ws.Range("N141").Value = -590224.5600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 251848.42
$ws.Range("J5").Value = 419064.78
$ws.Range("L5").Value = 1257194.34
$ws.Range("N5").Value = -1257418.34
$ws.Range("H6").Value = 41666800
$ws.Range("I6").Value = 41666800
$ws.Range("K6").Value = 125000400
$ws.Range("M6").Value = -125000287
$ws.Range("H11").Value = 582905.1
$ws.Range("I11").Value = 705998.4399999999
$ws.Range("J11").Value = 143286
$ws.Range("K11").Value = 2117995.32
$ws.Range("L11").Value = 429858
$ws.Range("M11").Value = -2117855.32
$ws.Range("N11").Value = -430138
$ws.Range("H12").Value = 85.5
$ws.Range("J12").Value = 85.875
$ws.Range("L12").Value = 257.625
$ws.Range("N12").Value = -603.625
$ws.Range("H22").Value = 987.8125
$ws.Range("I22").Value = 788.6
$ws.Range("J22").Value = 1319.8334
$ws.Range("K22").Value = 2365.8
$ws.Range("L22").Value = 3959.5002
$ws.Range("M22").Value = -2196.8
$ws.Range("N22").Value = -4297.5002
$ws.Range("H27").Value = 987.8125
$ws.Range("I27").Value = 788.6
$ws.Range("J27").Value = 1319.8334
$ws.Range("K27").Value = 2365.8
$ws.Range("L27").Value = 3959.5002
$ws.Range("M27").Value = -2263.8
$ws.Range("N27").Value = -4163.5002
$ws.Range("H32").Value = 872.75
$ws.Range("I32").Value = 400
$ws.Range("K32").Value = 1200
$ws.Range("M32").Value = -917
$ws.Range("H39").Value = 299
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
$ws.Range("H55").Value = 10121.652
$ws.Range("J55").Value = 11565.1
$ws.Range("L55").Value = 34695.3
$ws.Range("N55").Value = -35049.3
$ws.Range("H56").Value = 7982.5713
$ws.Range("I56").Value = 7982.5713
$ws.Range("K56").Value = 7982.5713
$ws.Range("M56").Value = -7452.5713
$ws.Range("H63").Value = 3059.5881
$ws.Range("J63").Value = 3000
$ws.Range("L63").Value = 9000
$ws.Range("N63").Value = -10498
$ws.Range("H66").Value = 3059.5881
$ws.Range("J66").Value = 3000
$ws.Range("L66").Value = 27000
$ws.Range("N66").Value = -34488
$ws.Range("H103").Value = 5072.75
$ws.Range("I103").Value = 5699.8
$ws.Range("K103").Value = 17099.4
$ws.Range("M103").Value = -16220.4
$ws.Range("H122").Value = 6288.3105
$ws.Range("I122").Value = 1407.5
$ws.Range("J122").Value = 8147.6665
$ws.Range("K122").Value = 12667.5
$ws.Range("L122").Value = 73328.9985
$ws.Range("M122").Value = -10217.5
$ws.Range("N122").Value = -78228.9985
$ws.Range("H129").Value = 1277
$ws.Range("I129").Value = 1045.1333
$ws.Range("J129").Value = 3016
$ws.Range("K129").Value = 3135.3999
$ws.Range("L129").Value = 9048
$ws.Range("M129").Value = 1864.6001
$ws.Range("N129").Value = -19048
$ws.Range("H131").Value = 6120.85
$ws.Range("I131").Value = 8410.299999999999
$ws.Range("K131").Value = 25230.9
$ws.Range("M131").Value = -20190.9
$ws.Range("H135").Value = 251848.42
$ws.Range("J135").Value = 419064.78
$ws.Range("L135").Value = 3771583.02
$ws.Range("N135").Value = -3776653.02

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H6").Value = 0
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H16").Value = 0
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").ClearContents()
$ws.Range("H102").Value = 3653.17
$ws.Range("I102").Value = 3905.6667
$ws.Range("J102").Value = 2232.875
$ws.Range("K102").Value = 3905.6667
$ws.Range("L102").Value = 2232.875
$ws.Range("M102").Value = -2283.6667
$ws.Range("N102").Value = -5476.875
$ws.Range("H113").Value = 39666.332
$ws.Range("J113").Value = 5000
$ws.Range("L113").Value = 5000
$ws.Range("N113").Value = -9340
$ws.Range("H126").Value = 10855.177
$ws.Range("I126").Value = 14960
$ws.Range("J126").Value = 4991.143
$ws.Range("K126").Value = 44880
$ws.Range("L126").Value = 14973.429
$ws.Range("M126").Value = -42410
$ws.Range("N126").Value = -19913.429
$ws.Range("H132").Value = 3761.2083
$ws.Range("I132").Value = 3856.4
$ws.Range("J132").Value = 2333.3333
$ws.Range("K132").Value = 11569.2
$ws.Range("L132").Value = 6999.999899999999
$ws.Range("M132").Value = -9039.200000000001
$ws.Range("N132").Value = -12059.9999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 696.5
$ws.Range("I22").Value = 696.5
$ws.Range("K22").Value = 696.5
$ws.Range("M22").Value = -401.5
$ws.Range("H25").Value = 6666.6665
$ws.Range("J25").Value = 7500
$ws.Range("L25").Value = 7500
$ws.Range("N25").Value = -7960
$ws.Range("H27").Value = 696.5
$ws.Range("I27").Value = 696.5
$ws.Range("K27").Value = 696.5
$ws.Range("M27").Value = -589.5
$ws.Range("H46").Value = 1385.5
$ws.Range("I46").Value = 919.1429000000001
$ws.Range("J46").Value = 2038.4
$ws.Range("K46").Value = 919.1429000000001
$ws.Range("L46").Value = 2038.4
$ws.Range("M46").Value = -731.1429000000001
$ws.Range("N46").Value = -2414.4
$ws.Range("H61").Value = 30671.455
$ws.Range("I61").Value = 2977.2
$ws.Range("J61").Value = 53750
$ws.Range("K61").Value = 2977.2
$ws.Range("L61").Value = 53750
$ws.Range("M61").Value = -2775.2
$ws.Range("N61").Value = -54154
$ws.Range("H68").Value = 4380.625
$ws.Range("I68").Value = 1817.3334
$ws.Range("J68").Value = 5918.6
$ws.Range("K68").Value = 1817.3334
$ws.Range("L68").Value = 5918.6
$ws.Range("M68").Value = -1068.3334
$ws.Range("N68").Value = -7416.6
$ws.Range("H71").Value = 4380.625
$ws.Range("I71").Value = 1817.3334
$ws.Range("J71").Value = 5918.6
$ws.Range("K71").Value = 9086.666999999999
$ws.Range("L71").Value = 29593
$ws.Range("M71").Value = -5342.666999999999
$ws.Range("N71").Value = -37081
$ws.Range("H93").Value = 5869.1875
$ws.Range("I93").Value = 6127.1333
$ws.Range("K93").Value = 6127.1333
$ws.Range("M93").Value = -4879.1333
$ws.Range("H100").Value = 5309.091
$ws.Range("I100").Value = 3080.2
$ws.Range("K100").Value = 3080.2
$ws.Range("M100").Value = -2539.2
$ws.Range("H113").Value = 30671.455
$ws.Range("I113").Value = 2977.2
$ws.Range("J113").Value = 53750
$ws.Range("K113").Value = 2977.2
$ws.Range("L113").Value = 53750
$ws.Range("M113").Value = -807.1999999999998
$ws.Range("N113").Value = -58090
$ws.Range("H122").Value = 8541.286
$ws.Range("I122").Value = 11539
$ws.Range("K122").Value = 34617
$ws.Range("M122").Value = -32167
$ws.Range("H132").Value = 935400.9
$ws.Range("I132").Value = 1658379.9
$ws.Range("K132").Value = 4975139.699999999
$ws.Range("M132").Value = -4972609.699999999
$ws.Range("H136").Value = 9480.166999999999
$ws.Range("I136").Value = 12183.167
$ws.Range("J136").Value = 8128.6665
$ws.Range("K136").Value = 36549.501
$ws.Range("L136").Value = 24385.9995
$ws.Range("M136").Value = -33999.501
$ws.Range("N136").Value = -29485.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 303821.78
$ws.Range("I62").Value = 303821.78
$ws.Range("K62").Value = 303821.78
$ws.Range("M62").Value = -303197.78
$ws.Range("H65").Value = 303821.78
$ws.Range("I65").Value = 303821.78
$ws.Range("K65").Value = 1519108.9
$ws.Range("M65").Value = -1515988.9
$ws.Range("H81").Value = 12868.091
$ws.Range("J81").Value = 5249.8335
$ws.Range("L81").Value = 10499.667
$ws.Range("N81").Value = -12621.667
$ws.Range("H84").Value = 12868.091
$ws.Range("J84").Value = 5249.8335
$ws.Range("L84").Value = 52498.335
$ws.Range("N84").Value = -63106.335
$ws.Range("H113").Value = 1972.1052
$ws.Range("I113").Value = 879.25
$ws.Range("J113").Value = 7800.6665
$ws.Range("K113").Value = 2637.75
$ws.Range("L113").Value = 23401.9995
$ws.Range("M113").Value = -467.75
$ws.Range("N113").Value = -27741.9995
$ws.Range("H126").Value = 27288.422
$ws.Range("I126").Value = 32299.066
$ws.Range("K126").Value = 96897.198
$ws.Range("M126").Value = -94427.198
$ws.Range("H132").Value = 11783.102
$ws.Range("I132").Value = 13384.45
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 40153.35000000001
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -37623.35000000001
$ws.Range("N132").Value = -19058
